$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "PEKeK595"
$ws.Range("B2").Value = 23090411
$ws.Range("C2").Value = "wrjxavd33"
$ws.Range("D2").Value = "PmbY`$5!7"
$ws.Range("F2").Value = "lmFrXCvp"
$ws.Range("G2").Value = "Kssc"
